$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.106.81'
$ws.Range('E2').Value = '  +5.71%  '
$ws.Range('D3').Value = '1.919.97'
$ws.Range('E3').Value = '  +2.63%  '
$ws.Range('E4').Value = '  -0.56%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '330.14'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  +4.66%  '
$ws.Range('E6').Value = '  -0.54%  '
$origStyle = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5196'
$ws.Range('D7').Style = $origStyle
$ws.Range('E7').Value = '  +1.99%  '
$origStyle = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4082'
$ws.Range('D8').Style = $origStyle
$ws.Range('E8').Value = '  +4.71%  '
$origStyle = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08506'
$ws.Range('D9').Style = $origStyle
$ws.Range('E9').Value = '  +1.97%  '
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '43.26'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  +3.64%  '
$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.127'
$ws.Range('D11').Style = $origStyle
$ws.Range('E11').Value = '  +2.13%  '
$origStyle = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.39'
$ws.Range('D12').Style = $origStyle
$ws.Range('E12').Value = '  +9.69%  '
$origStyle = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.404'
$ws.Range('D13').Style = $origStyle
$ws.Range('E13').Value = '  +3.12%  '
$ws.Range('D14').Value = '1.922.87'
$ws.Range('E14').Value = '  +2.81%  '
$origStyle = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.403'
$ws.Range('D15').Style = $origStyle
$ws.Range('E15').Value = '  +1.72%  '
$ws.Range('E16').Value = '  -0.58%  '
$origStyle = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '95.54'
$ws.Range('D17').Style = $origStyle
$ws.Range('E17').Value = '  +4.90%  '
$ws.Range('E18').Value = '  +1.00%  '
$origStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06714'
$ws.Range('D19').Style = $origStyle
$ws.Range('E19').Value = '  -0.17%  '
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.25'
$ws.Range('D20').Style = $origStyle
$ws.Range('E20').Value = '  +3.22%  '
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('D21').Style = $origStyle
$ws.Range('E21').Value = '  -0.57%  '
$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.020'
$ws.Range('D22').Style = $origStyle
$ws.Range('E22').Value = '  +1.90%  '
$ws.Range('D23').Value = '30.115.02'
$ws.Range('E23').Value = '  +5.66%  '
$origStyle = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.33'
$ws.Range('D24').Style = $origStyle
$ws.Range('E24').Value = '  +1.95%  '
$ws.Range('E25').Value = '  +0.37%  '
$ws.Range('D26').Value = '2.142.04'
$ws.Range('E26').Value = '  +2.80%  '
$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.09'
$ws.Range('D27').Style = $origStyle
$ws.Range('E27').Value = '  +2.26%  '
$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '160.06'
$ws.Range('D28').Style = $origStyle
$ws.Range('E28').Value = '  -0.40%  '
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.448'
$ws.Range('D29').Style = $origStyle
$ws.Range('E29').Value = '  +1.74%  '
$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '129.01'
$ws.Range('D30').Style = $origStyle
$ws.Range('E30').Value = '  +2.16%  '
$origStyle = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.074'
$ws.Range('D31').Style = $origStyle
$ws.Range('E31').Value = '  +3.69%  '
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1054'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  +1.35%  '
$origStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.084'
$ws.Range('D33').Style = $origStyle
$ws.Range('E33').Value = '  +5.99%  '
$origStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.637'
$ws.Range('D34').Style = $origStyle
$ws.Range('E34').Value = '  +0.77%  '
$origStyle = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02493'
$ws.Range('D35').Style = $origStyle
$ws.Range('E35').Value = '  +1.78%  '
$origStyle = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06606'
$ws.Range('D36').Style = $origStyle
$ws.Range('E36').Value = '  +0.90%  '
$origStyle = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2205'
$ws.Range('D37').Style = $origStyle
$ws.Range('E37').Value = '  +1.99%  '
$origStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.202'
$ws.Range('D38').Style = $origStyle
$ws.Range('E38').Value = '  +3.72%  '
$origStyle = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.232'
$ws.Range('D39').Style = $origStyle
$ws.Range('E39').Value = '  +4.48%  '
$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.916'
$ws.Range('D40').Style = $origStyle
$ws.Range('E40').Value = '  +0.35%  '
$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6516'
$ws.Range('D41').Style = $origStyle
$ws.Range('E41').Value = '  +2.37%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$origStyle = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.248'
$ws.Range('D42').Style = $origStyle
$ws.Range('E42').Value = '  +1.14%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$origStyle = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.63'
$ws.Range('D43').Style = $origStyle
$ws.Range('E43').Value = '  +4.86%  '
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6163'
$ws.Range('D44').Style = $origStyle
$ws.Range('E44').Value = '  +2.72%  '
$origStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.22'
$ws.Range('D45').Style = $origStyle
$ws.Range('E45').Value = '  +1.61%  '
$origStyle = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.772'
$ws.Range('D46').Style = $origStyle
$ws.Range('E46').Value = '  +2.20%  '
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.083'
$ws.Range('D47').Style = $origStyle
$ws.Range('E47').Value = '  +4.22%  '
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.246'
$ws.Range('D48').Style = $origStyle
$ws.Range('E48').Value = '  +2.49%  '
$origStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '124.38'
$ws.Range('D49').Style = $origStyle
$ws.Range('E49').Value = '  +2.07%  '
$origStyle = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.162'
$ws.Range('D50').Style = $origStyle
$ws.Range('E50').Value = '  +3.09%  '
$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.73'
$ws.Range('D51').Style = $origStyle
$ws.Range('E51').Value = '  +4.58%  '
